$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert two new rows before row 15. This pushes the old row 15 (totals)
#    and row 16 (footer) down to rows 17 and 18, carrying their values,
#    styles and merged cells with them.
$ws.Rows("15:15").Insert()
$ws.Rows("15:15").Insert()

# 2. Populate the two new sale rows (15 and 16) with their values first -
#    text-like numeric strings get a leading apostrophe so they are stored
#    as text (matching the source data) instead of being coerced to numbers.
$ws.Range("A15").Value = 9
$ws.Range("C15").Value = "فرشه شعر اطفال الجو"
$ws.Range("H15").Value = "4:0"
$ws.Range("L15").Value = "'0"
$ws.Range("N15").Value = "'25.00"
$ws.Range("P15").Value = "'25.0000"
$ws.Range("Q15").Value = "1:0"

$ws.Range("A16").Value = 10
$ws.Range("C16").Value = "مخمريه العود الملكي"
$ws.Range("H16").Value = "3:0"
$ws.Range("L16").Value = "'0"
$ws.Range("N16").Value = "'35.00"
$ws.Range("P16").Value = "'35.0000"
$ws.Range("Q16").Value = "1:0"

# 3. Copy the formatting of an existing item row (row 7) onto the two new
#    rows so they match the look of the other sale rows.
$ws.Range("A7:Q7").Copy()
$ws.Range("A15:Q15").PasteSpecial(-4122)
$ws.Range("A7:Q7").Copy()
$ws.Range("A16:Q16").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 4. Row heights for the new rows (row 17, the shifted totals row, ends up
#    at 25.5 to match the source workbook after the insert).
$ws.Rows(15).RowHeight = 24.75
$ws.Rows(16).RowHeight = 25.5
$ws.Rows(17).RowHeight = 25.5

# 5. Re-create the merged cells for the new rows (mirrors the merge layout
#    used by every other item row).
$ws.Range("A15:B15").Merge()
$ws.Range("C15:G15").Merge()
$ws.Range("H15:K15").Merge()
$ws.Range("L15:M15").Merge()
$ws.Range("N15:O15").Merge()

$ws.Range("A16:B16").Merge()
$ws.Range("C16:G16").Merge()
$ws.Range("H16:K16").Merge()
$ws.Range("L16:M16").Merge()
$ws.Range("N16:O16").Merge()

# 6. Update the running total (now on row 17) to include the two new items.
$ws.Range("P17").Value = 328.89999999999998

# 7. Update the printed timestamp (now on row 18) to the new generation time.
$ws.Range("A18").Value = "Saturday, 13 September, 2025 11:00 AM"
